$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new parametrisation columns: age (D) and gender (E)
# Order of writes matters for shared-string table ordering, so follow the
# row-by-row / left-to-right layout of the final sheet.
$ws.Range("D1").Value = "age"
$ws.Range("E1").Value = "gender"

$ws.Range("D2").Value = 22
$ws.Range("E2").Value = "M"

$ws.Range("D3").Value = 23
$ws.Range("E3").Value = "F"

# Update existing username value in A3 (was "Manish1", now "TTT")
$ws.Range("A3").Value = "TTT"

# Update selection to match target state
$ws.Range("A3").Select()
